$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.128.45'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '3.659.89'
$ws.Range("E3").Value = '  -0.36%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.01'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.17%  '

$ws.Range("E7").Value = '  -0.41%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.696'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.96%  '

$ws.Range("E10").Value = '  -5.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.59'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000270'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -5.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.22'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.64%  '

$ws.Range("D14").Value = '4.251.33'
$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("D15").Value = '3.660.45'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("E16").Value = '  +0.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.85'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.11%  '

$ws.Range("D18").Value = '67.929.63'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("E20").Value = '  -1.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '404.01'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.42'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '87.95'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.85'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.48'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.58%  '

$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.32'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.85'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.15'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '68.01'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.22'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '44.08'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.04%  '

$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '605.99'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.58%  '

$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.390'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.38%  '

$ws.Range("D40").Value = '0.0₃0765'
$ws.Range("E40").Value = '  -12.62%  '

$ws.Range("E41").Value = '  +0.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0423'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.10%  '

$ws.Range("E44").Value = '  -8.56%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.23'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.22%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.135'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.74%  '

$ws.Range("D47").Value = '2.762.75'
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.87'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.31'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("E50").Value = '  -4.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.53'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -10.21%  '

